$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feature: mark algorithms containing the pile-of-poo emoji as "difficult".
# Append the emoji to the two algorithm strings that are now flagged difficult.
$ws.Range("B3").Value = "U' [U perm] 💩"
$ws.Range("D3").Value = "[U perm] 💩"
$ws.Range("C4").Value = "[U perm] 💩"

# The previously-empty placeholder string in E3 is no longer needed now that
# it isn't rendering a forced blank string - clear it to a real empty cell.
$ws.Range("E3").ClearContents()

# Left-align columns A, B and D (the row-label / algorithm columns).
$ws.Range("A:A").HorizontalAlignment = -4131
$ws.Range("B:B").HorizontalAlignment = -4131
$ws.Range("D:D").HorizontalAlignment = -4131

# Slightly taller rows to accommodate the emoji glyph.
$ws.Range("A1:E5").RowHeight = 18.75
